# Generate Report for Handback
# The 6869bf9d-... file has now been handed back in sync with en-US, so refresh
# the status / handback datetime / error-detail columns across the report sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the 6869bf9d-... file ---
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E3").Value = "Handed back: in sync with en-US"
$ovw.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row 3 is the 6869bf9d-... file ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-09-02 20:56:32"
$zhcn.Range("P3").Value = ""

# --- de-de sheet: row 3 is the 6869bf9d-... file ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-09-02 20:56:40"
$dede.Range("P3").Value = ""
